$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; this shifts the old D:K data to F:M
$ws.Range("D5:E5").EntireColumn.Insert()

# Copy number formats (date format / number format) from column F (the old column D,
# now shifted two columns right) into the brand-new D:E columns, per contiguous data block,
# so blank "spacer" rows that never had D:K cells stay untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns (D, E) with the new quarters of data
$ws.Range("D7").Value = 43464
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2656800
$ws.Range("E8").Value = 2697600
$ws.Range("D9").Value = 2544900
$ws.Range("E9").Value = 2527900
$ws.Range("D10").Value = 111900
$ws.Range("E10").Value = 169700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 2600
$ws.Range("E14").Value = 300
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 2633200
$ws.Range("E17").Value = 2612300
$ws.Range("D18").Value = 23600
$ws.Range("E18").Value = 85300
$ws.Range("D20").Value = -15900
$ws.Range("E20").Value = 10300
$ws.Range("D21").Value = 75900
$ws.Range("E21").Value = 167600
$ws.Range("D22").Value = 36900
$ws.Range("E22").Value = 35300
$ws.Range("D23").Value = -29200
$ws.Range("E23").Value = 60300
$ws.Range("D24").Value = -47300
$ws.Range("E24").Value = 30800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 18200
$ws.Range("E26").Value = 29500
$ws.Range("D27").Value = 19100
$ws.Range("E27").Value = 29300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -26400
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 15900
$ws.Range("E32").Value = -10300
$ws.Range("D33").Value = -7300
$ws.Range("E33").Value = 29300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -7300
$ws.Range("E35").Value = 29300
$ws.Range("D38").Value = 43464
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 338400
$ws.Range("E41").Value = 401300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 601400
$ws.Range("E43").Value = 564800
$ws.Range("D44").Value = 1159500
$ws.Range("E44").Value = 1183000
$ws.Range("D45").Value = 135400
$ws.Range("E45").Value = 133700
$ws.Range("D46").Value = 2234700
$ws.Range("E46").Value = 2282700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2161700
$ws.Range("E48").Value = 2137800
$ws.Range("D49").Value = 1513900
$ws.Range("E49").Value = 1552700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 21000
$ws.Range("E52").Value = 3200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 5931200
$ws.Range("E54").Value = 5976400
$ws.Range("D57").Value = 744100
$ws.Range("E57").Value = 740700
$ws.Range("D58").Value = 30400
$ws.Range("E58").Value = 24000
$ws.Range("D59").Value = 521700
$ws.Range("E59").Value = 558400
$ws.Range("D60").Value = 1296200
$ws.Range("E60").Value = 1323100
$ws.Range("D61").Value = 2295200
$ws.Range("E61").Value = 2302200
$ws.Range("D62").Value = 320200
$ws.Range("E62").Value = 283600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3921400
$ws.Range("E66").Value = 3918200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 421900
$ws.Range("E72").Value = 429200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 2009800
$ws.Range("E76").Value = 2058300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43464
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -7300
$ws.Range("E81").Value = 29300
$ws.Range("D83").Value = 68200
$ws.Range("E83").Value = 72000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 66400
$ws.Range("E89").Value = 121300
$ws.Range("D91").Value = -116800
$ws.Range("E91").Value = -76700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -109700
$ws.Range("E94").Value = -75200
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -10400
$ws.Range("E100").Value = -304300
$ws.Range("D101").Value = -500
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -54300
$ws.Range("E102").Value = -258200

# Data corrections that came bundled with this update (not simple column shifts)
$ws.Range("I91").Value = -84200
$ws.Range("J91").Value = -59700
